# The workbook gained one new weekly data row: a row is inserted at row 62
# (pushing the existing rows 62-100 down to 63-101 unchanged), and the newly
# inserted row 62 repeats the same market/category/variety/price data as the
# row that used to be there, but dated one reporting period later.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 62; everything below shifts down.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new week's record.
$ws.Range("A62").Value = 4
$ws.Range("B62").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C62").Value = "Los Lagos"
$ws.Range("D62").Value = 44827
$ws.Range("E62").Value = 10
$ws.Range("F62").Value = 100112031
$ws.Range("G62").Value = "Poroto verde"
$ws.Range("H62").Value = "Magnum"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 35
$ws.Range("K62").Value = 35000
$ws.Range("L62").Value = 35000
$ws.Range("M62").Value = 35000
$ws.Range("N62").Value = "$/malla 25 kilos"
$ws.Range("O62").Value = "Perú"
$ws.Range("P62").Value = 1400
$ws.Range("Q62").Value = 25
$ws.Range("R62").Value = "Hortaliza"
